# Rail pathless OD pairs are now derived to the roadway network during
# FreightNetwork construction. This shifts several "Roadway" column values
# (C, E, G, I, K, M) on the global_results sheet - update them to their
# newly-recomputed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("global_results")

$ws.Range("C4").Value = 0.005665227045420121
$ws.Range("E4").Value = 0.007470363524626595
$ws.Range("G4").Value = 0.005363387304790163
$ws.Range("I4").Value = 0.005682562833820645
$ws.Range("K4").Value = 0.007322221694108118
$ws.Range("M4").Value = 0.005363387304790163
$ws.Range("C6").Value = 0.06566522704542013
$ws.Range("E6").Value = 0.0674703635246266
$ws.Range("G6").Value = 0.06536338730479017
$ws.Range("I6").Value = 0.06568256283382065
$ws.Range("K6").Value = 0.06732222169410812
$ws.Range("M6").Value = 0.06536338730479017
$ws.Range("C8").Value = 409093579.8324731
$ws.Range("E8").Value = 342945418.7755197
$ws.Range("G8").Value = 439376480.5473356
$ws.Range("I8").Value = 409501679.8324731
$ws.Range("K8").Value = 349006724.5486711
$ws.Range("M8").Value = 439376480.5473356
$ws.Range("C9").Value = 122009715530.0625
$ws.Range("E9").Value = 90318659924.64212
$ws.Range("G9").Value = 133245739070.795
$ws.Range("I9").Value = 122175050590.0625
$ws.Range("K9").Value = 92899577976.22893
$ws.Range("M9").Value = 133245739070.795
$ws.Range("C10").Value = 8011795672.02868
$ws.Range("E10").Value = 6093832818.172728
$ws.Range("G10").Value = 8709392849.597383
$ws.Range("I10").Value = 8024770437.107
$ws.Range("K10").Value = 6254205983.804768
$ws.Range("M10").Value = 8709392849.597383
$ws.Range("C15").Value = 298.2440242157465
$ws.Range("E15").Value = 263.3616166885192
$ws.Range("G15").Value = 303.2609731517932
$ws.Range("I15").Value = 298.3505480125119
$ws.Range("K15").Value = 266.1827736882861
$ws.Range("M15").Value = 303.2609731517932
$ws.Range("I16").Value = 33131.8
$ws.Range("K16").Value = 33131.8
$ws.Range("M16").Value = 33911.8
$ws.Range("C17").Value = 3704808.444166987
$ws.Range("E17").Value = 2742513.844089847
$ws.Range("G17").Value = 3929185.094002529
$ws.Range("I17").Value = 3687546.423377617
$ws.Range("K17").Value = 2803939.960286761
$ws.Range("M17").Value = 3929185.094002529

# The blank (empty-but-string-typed) cells M7, M11, M12 and M13 otherwise get
# normalized by the round-trip into a visible shared string on save; force
# them back to a true empty value so they remain blank, as in the source file.
$ws.Range("M7").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("M13").Value = ""
